$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(7, 7).Value = 1.36
$ws.Cells.Item(7, 8).Value = 4.5
$ws.Cells.Item(7, 9).Value = 9.5
$ws.Cells.Item(7, 10).Value = 1.91
$ws.Cells.Item(7, 11).Value = 2.2
$ws.Cells.Item(7, 12).Value = 9.5
$ws.Cells.Item(7, 15).Value = 1.36
$ws.Cells.Item(7, 16).Value = 3
$ws.Cells.Item(7, 17).Value = 2.15
$ws.Cells.Item(7, 18).Value = 1.67
$ws.Cells.Item(7, 19).Value = 3.15
$ws.Cells.Item(7, 20).Value = 1.37
$ws.Cells.Item(7, 21).Value = 4
$ws.Cells.Item(7, 22).Value = 1.22
$ws.Cells.Item(7, 25).Value = 2.63
$ws.Cells.Item(7, 26).Value = 1.44
$ws.Cells.Item(7, 28).Value = 5.5
$ws.Cells.Item(7, 29).Value = 9.5
$ws.Cells.Item(7, 30).Value = 8
$ws.Cells.Item(7, 33).Value = 8
$ws.Cells.Item(7, 34).Value = 9
$ws.Cells.Item(7, 35).Value = 29
$ws.Cells.Item(7, 36).Value = 126
$ws.Cells.Item(7, 38).Value = 15
$ws.Cells.Item(7, 39).Value = 41
$ws.Cells.Item(7, 40).Value = 29
$ws.Cells.Item(7, 41).Value = 126
$ws.Cells.Item(7, 42).Value = 81
$ws.Cells.Item(7, 43).Value = 81
$ws.Cells.Item(7, 44).Value = 1.61
$ws.Cells.Item(7, 45).Value = 2.33
$ws.Cells.Item(8, 8).Value = 2.75
$ws.Cells.Item(8, 10).Value = 4.5
$ws.Cells.Item(8, 13).Value = 1.17
$ws.Cells.Item(8, 14).Value = 5
$ws.Cells.Item(8, 15).Value = 1.73
$ws.Cells.Item(8, 16).Value = 2
$ws.Cells.Item(8, 17).Value = 3.5
$ws.Cells.Item(8, 18).Value = 1.3
$ws.Cells.Item(8, 21).Value = 8
$ws.Cells.Item(8, 22).Value = 1.08
$ws.Cells.Item(8, 23).Value = 1.78
$ws.Cells.Item(8, 24).Value = 2.03
$ws.Cells.Item(8, 33).Value = 4.75
$ws.Cells.Item(8, 38).Value = 5
$ws.Cells.Item(9, 8).Value = 2.9
$ws.Cells.Item(9, 13).Value = 1.17
$ws.Cells.Item(9, 14).Value = 5
$ws.Cells.Item(9, 33).Value = 5
$ws.Cells.Item(9, 34).Value = 6.5
$ws.Cells.Item(24, 7).Value = 2.38
$ws.Cells.Item(24, 8).Value = 3.5
$ws.Cells.Item(24, 9).Value = 2.63
$ws.Cells.Item(24, 10).Value = 3.2
$ws.Cells.Item(24, 12).Value = 3.5
$ws.Cells.Item(24, 13).Value = 1.06
$ws.Cells.Item(24, 14).Value = 10
$ws.Cells.Item(24, 15).Value = 1.33
$ws.Cells.Item(24, 16).Value = 3.25
$ws.Cells.Item(24, 17).Value = 2.05
$ws.Cells.Item(24, 18).Value = 1.8
$ws.Cells.Item(24, 19).Value = 2.85
$ws.Cells.Item(24, 20).Value = 1.41
$ws.Cells.Item(24, 23).Value = 1.44
$ws.Cells.Item(24, 24).Value = 2.63
$ws.Cells.Item(24, 25).Value = 1.83
$ws.Cells.Item(24, 26).Value = 1.83
$ws.Cells.Item(24, 30).Value = 23
$ws.Cells.Item(24, 31).Value = 21
$ws.Cells.Item(24, 33).Value = 10
$ws.Cells.Item(24, 38).Value = 8.5
$ws.Cells.Item(24, 39).Value = 13
$ws.Cells.Item(24, 40).Value = 10
$ws.Cells.Item(24, 44).Value = 1.53
$ws.Cells.Item(24, 45).Value = 2.47
$ws.Cells.Item(25, 7).Value = 1.85
$ws.Cells.Item(25, 8).Value = 3.2
$ws.Cells.Item(25, 9).Value = 4.75
$ws.Cells.Item(25, 10).Value = 2.63
$ws.Cells.Item(25, 11).Value = 1.91
$ws.Cells.Item(25, 12).Value = 5.5
$ws.Cells.Item(25, 13).Value = 1.11
$ws.Cells.Item(25, 14).Value = 6.5
$ws.Cells.Item(25, 15).Value = 1.53
$ws.Cells.Item(25, 16).Value = 2.38
$ws.Cells.Item(25, 17).Value = 2.7
$ws.Cells.Item(25, 18).Value = 1.44
$ws.Cells.Item(25, 19).Value = 4.3
$ws.Cells.Item(25, 20).Value = 1.21
$ws.Cells.Item(25, 23).Value = 1.62
$ws.Cells.Item(25, 24).Value = 2.2
$ws.Cells.Item(25, 30).Value = 15
$ws.Cells.Item(25, 31).Value = 21
$ws.Cells.Item(25, 33).Value = 6
$ws.Cells.Item(25, 36).Value = 101
$ws.Cells.Item(25, 38).Value = 9
$ws.Cells.Item(25, 39).Value = 21
$ws.Cells.Item(25, 42).Value = 41
$ws.Cells.Item(25, 44).Value = 2.05
$ws.Cells.Item(25, 45).Value = 1.8
$ws.Cells.Item(26, 7).Value = 2.25
$ws.Cells.Item(26, 8).Value = 2.88
$ws.Cells.Item(26, 9).Value = 3.7
$ws.Cells.Item(26, 11).Value = 1.8
$ws.Cells.Item(26, 13).Value = 1.17
$ws.Cells.Item(26, 14).Value = 5
$ws.Cells.Item(26, 15).Value = 1.67
$ws.Cells.Item(26, 16).Value = 2.1
$ws.Cells.Item(26, 28).Value = 9
$ws.Cells.Item(26, 31).Value = 26
$ws.Cells.Item(26, 33).Value = 5
$ws.Cells.Item(26, 34).Value = 6
$ws.Cells.Item(26, 39).Value = 17
$ws.Cells.Item(35, 7).Value = 1.73
$ws.Cells.Item(35, 8).Value = 3.7
$ws.Cells.Item(35, 9).Value = 4.5
$ws.Cells.Item(35, 10).Value = 2.38
$ws.Cells.Item(35, 17).Value = 1.85
$ws.Cells.Item(35, 18).Value = 1.95
$ws.Cells.Item(35, 23).Value = 1.36
$ws.Cells.Item(35, 24).Value = 3
$ws.Cells.Item(54, 7).Value = 2.1
$ws.Cells.Item(54, 8).Value = 3.7
$ws.Cells.Item(54, 9).Value = 3.2
$ws.Cells.Item(54, 10).Value = 2.63
$ws.Cells.Item(54, 15).Value = 1.17
$ws.Cells.Item(54, 16).Value = 5
$ws.Cells.Item(54, 17).Value = 1.57
$ws.Cells.Item(54, 18).Value = 2.35
$ws.Cells.Item(54, 19).Value = 1.98
$ws.Cells.Item(54, 20).Value = 1.88
$ws.Cells.Item(54, 21).Value = 2.38
$ws.Cells.Item(54, 22).Value = 1.53
$ws.Cells.Item(54, 27).Value = 11
$ws.Cells.Item(54, 30).Value = 21
$ws.Cells.Item(54, 33).Value = 15
$ws.Cells.Item(55, 7).Value = 2.05
$ws.Cells.Item(55, 9).Value = 3.75
$ws.Cells.Item(55, 25).Value = 2
$ws.Cells.Item(55, 26).Value = 1.75
$ws.Cells.Item(56, 8).Value = 5.1
$ws.Cells.Item(56, 9).Value = 8.75
$ws.Cells.Item(56, 11).Value = 2.55
$ws.Cells.Item(56, 12).Value = 7.1
$ws.Cells.Item(56, 16).Value = 5
$ws.Cells.Item(56, 17).Value = 1.47
$ws.Cells.Item(56, 18).Value = 2.32
$ws.Cells.Item(56, 21).Value = 2.12
$ws.Cells.Item(56, 22).Value = 1.57
$ws.Cells.Item(56, 25).Value = 1.82
$ws.Cells.Item(56, 27).Value = 8.25
$ws.Cells.Item(56, 28).Value = 6.9
$ws.Cells.Item(56, 33).Value = 16
$ws.Cells.Item(56, 34).Value = 10.5
$ws.Cells.Item(56, 37).Value = 600
$ws.Cells.Item(56, 39).Value = 70
$ws.Cells.Item(56, 42).Value = 100
$ws.Cells.Item(56, 43).Value = 75
$ws.Cells.Item(57, 7).Value = 1.5
$ws.Cells.Item(57, 8).Value = 3.4
$ws.Cells.Item(57, 9).Value = 6.25
$ws.Cells.Item(57, 10).Value = 2.2
$ws.Cells.Item(57, 12).Value = 7
$ws.Cells.Item(57, 21).Value = 4.33
$ws.Cells.Item(57, 22).Value = 1.2
$ws.Cells.Item(57, 23).Value = 1.5
$ws.Cells.Item(57, 24).Value = 2.5
$ws.Cells.Item(57, 28).Value = 6
$ws.Cells.Item(57, 30).Value = 10
$ws.Cells.Item(57, 31).Value = 15
$ws.Cells.Item(57, 34).Value = 7.5
$ws.Cells.Item(57, 36).Value = 101
$ws.Cells.Item(57, 38).Value = 13
$ws.Cells.Item(57, 39).Value = 34
$ws.Cells.Item(57, 41).Value = 81
$ws.Cells.Item(57, 44).Value = 1.75
$ws.Cells.Item(57, 45).Value = 2.05
$ws.Cells.Item(58, 15).Value = 1.57
$ws.Cells.Item(58, 16).Value = 2.25
$ws.Cells.Item(58, 21).Value = 6
$ws.Cells.Item(58, 22).Value = 1.13
$ws.Cells.Item(67, 7).Value = 1.91
$ws.Cells.Item(67, 8).Value = 2.9
$ws.Cells.Item(67, 15).Value = 1.44
$ws.Cells.Item(67, 16).Value = 2.63
$ws.Cells.Item(67, 17).Value = 2.4
$ws.Cells.Item(67, 18).Value = 1.53
$ws.Cells.Item(67, 21).Value = 4.5
$ws.Cells.Item(67, 22).Value = 1.18
$ws.Cells.Item(67, 23).Value = 1.53
$ws.Cells.Item(67, 24).Value = 2.38
$ws.Cells.Item(67, 31).Value = 19
$ws.Cells.Item(67, 33).Value = 6.5
$ws.Cells.Item(67, 37).Value = 1000
$ws.Cells.Item(67, 38).Value = 10
$ws.Cells.Item(67, 39).Value = 21
$ws.Cells.Item(67, 44).Value = 1.78
$ws.Cells.Item(67, 45).Value = 2.03
